# Updated CVDs for the month
$wb = $excel.ActiveWorkbook

# Sheet: Downers Grove Aerospace Illino
$ws3 = $wb.Worksheets.Item("Downers Grove Aerospace Illino")
$ws3.Range("E2").Value = 0.0732
$ws3.Range("E3").Value = 0.0732
$ws3.Range("E4").Value = 0.0732
$ws3.Range("O4").Value = 0.0102
$ws3.Range("P4").Value = 0
$ws3.Range("Q4").Value = 0
$ws3.Range("R4").Value = 0
$ws3.Range("S4").Value = 0
$ws3.Range("T4").Value = 0
$ws3.Range("U4").Value = 0
$ws3.Range("V4").Value = 0
$ws3.Range("W4").Value = 0
$ws3.Range("E5").Value = 0.272727272727273
$ws3.Range("E6").Value = 0.272727272727273
$ws3.Range("E7").Value = 0.272727272727273
$ws3.Range("O7").Value = 0
$ws3.Range("P7").Value = 0.272727272727273
$ws3.Range("Q7").Value = 0.272727272727273
$ws3.Range("R7").Value = 0.272727272727273
$ws3.Range("S7").Value = 0.272727272727273
$ws3.Range("T7").Value = 0.272727272727273
$ws3.Range("U7").Value = 0.272727272727273
$ws3.Range("V7").Value = 0.272727272727273
$ws3.Range("W7").Value = 0.272727272727273
$ws3.Range("E8").Value = 0.1248
$ws3.Range("E9").Value = 0.1248
$ws3.Range("E10").Value = 0.1248
$ws3.Range("G10").Value = 0.0228
$ws3.Range("H10").Value = 0.0092
$ws3.Range("I10").Value = 0.018
$ws3.Range("J10").Value = 0.0501
$ws3.Range("M10").Value = 0.0231
$ws3.Range("N10").Value = 0.051
$ws3.Range("O10").Value = 0.0239
$ws3.Range("P10").Value = 0.017825
$ws3.Range("Q10").Value = 0.017825
$ws3.Range("R10").Value = 0.053475
$ws3.Range("S10").Value = 0.017825
$ws3.Range("T10").Value = 0.017825
$ws3.Range("U10").Value = 0.017825
$ws3.Range("V10").Value = 0.053475
$ws3.Range("W10").Value = 0.2139

# Sheet: Monticello Indiana
$ws7 = $wb.Worksheets.Item("Monticello Indiana")
$ws7.Range("E2").Value = 0.0713
$ws7.Range("E3").Value = 0.0713
$ws7.Range("E4").Value = 0.0713
$ws7.Range("O4").Value = 0
$ws7.Range("P4").Value = 0
$ws7.Range("Q4").Value = 0
$ws7.Range("R4").Value = 0
$ws7.Range("S4").Value = 0
$ws7.Range("T4").Value = 0
$ws7.Range("U4").Value = 0
$ws7.Range("V4").Value = 0
$ws7.Range("W4").Value = 0
$ws7.Range("O7").Value = ""
$ws7.Range("E8").Value = 0.1003
$ws7.Range("E9").Value = 0.1003
$ws7.Range("E10").Value = 0.1003
$ws7.Range("I10").Value = 0.0375
$ws7.Range("J10").Value = 0.039
$ws7.Range("O10").Value = 0.012
$ws7.Range("P10").Value = 0.014325
$ws7.Range("Q10").Value = 0.014325
$ws7.Range("R10").Value = 0.042975
$ws7.Range("S10").Value = 0.014325
$ws7.Range("T10").Value = 0.014325
$ws7.Range("U10").Value = 0.014325
$ws7.Range("V10").Value = 0.042975
$ws7.Range("W10").Value = 0.1719

# Sheet: Silvestre Terrazas, Chihuahua 
$ws9 = $wb.Worksheets.Item("Silvestre Terrazas, Chihuahua ")
$ws9.Range("E2").Value = 0.0323
$ws9.Range("E3").Value = 0.0323
$ws9.Range("E4").Value = 0.0323
$ws9.Range("M4").Value = 0.0323
$ws9.Range("N4").Value = 0.0323
$ws9.Range("O4").Value = 0
$ws9.Range("P4").Value = 0
$ws9.Range("Q4").Value = 0
$ws9.Range("R4").Value = 0
$ws9.Range("S4").Value = 0
$ws9.Range("T4").Value = 0
$ws9.Range("U4").Value = 0
$ws9.Range("V4").Value = 0
$ws9.Range("W4").Value = 0
$ws9.Range("E5").Value = 0.375
$ws9.Range("E6").Value = 0.375
$ws9.Range("E7").Value = 0.375
$ws9.Range("O7").Value = 0
$ws9.Range("P7").Value = 0.375
$ws9.Range("Q7").Value = 0.375
$ws9.Range("R7").Value = 0.375
$ws9.Range("S7").Value = 0.375
$ws9.Range("T7").Value = 0.375
$ws9.Range("U7").Value = 0.375
$ws9.Range("V7").Value = 0.375
$ws9.Range("W7").Value = 0.375
$ws9.Range("E8").Value = 0.151
$ws9.Range("E9").Value = 0.151
$ws9.Range("E10").Value = 0.151
$ws9.Range("G10").Value = 0.0274
$ws9.Range("H10").Value = 0.0233
$ws9.Range("J10").Value = 0.049
$ws9.Range("K10").Value = 0.0116
$ws9.Range("L10").Value = 0.0211
$ws9.Range("M10").Value = 0.0095
$ws9.Range("N10").Value = 0.042
$ws9.Range("O10").Value = 0.0508
$ws9.Range("P10").Value = 0.021575
$ws9.Range("Q10").Value = 0.021575
$ws9.Range("R10").Value = 0.064725
$ws9.Range("S10").Value = 0.021575
$ws9.Range("T10").Value = 0.021575
$ws9.Range("U10").Value = 0.021575
$ws9.Range("V10").Value = 0.064725
$ws9.Range("W10").Value = 0.2589

# Sheet: Simi Valley California
$ws10 = $wb.Worksheets.Item("Simi Valley California")
$ws10.Range("E2").Value = 0.078
$ws10.Range("E3").Value = 0.078
$ws10.Range("E4").Value = 0.078
$ws10.Range("G4").Value = 0.0196
$ws10.Range("H4").Value = 0.0204
$ws10.Range("I4").Value = 0.02
$ws10.Range("J4").Value = 0.06
$ws10.Range("K4").Value = 0.0196
$ws10.Range("N4").Value = 0.0193
$ws10.Range("O4").Value = 0
$ws10.Range("P4").Value = 0
$ws10.Range("Q4").Value = 0
$ws10.Range("R4").Value = 0
$ws10.Range("S4").Value = 0
$ws10.Range("T4").Value = 0
$ws10.Range("U4").Value = 0
$ws10.Range("V4").Value = 0
$ws10.Range("W4").Value = 0
$ws10.Range("O7").Value = ""
$ws10.Range("E8").Value = 0.0525
$ws10.Range("E9").Value = 0.0525
$ws10.Range("E10").Value = 0.0525
$ws10.Range("G10").Value = 0.0055
$ws10.Range("H10").Value = 0.0055
$ws10.Range("I10").Value = 0.0054
$ws10.Range("J10").Value = 0.0163
$ws10.Range("M10").Value = 0.0149
$ws10.Range("N10").Value = 0.0206
$ws10.Range("O10").Value = 0.0151
$ws10.Range("P10").Value = 0.0075
$ws10.Range("Q10").Value = 0.0075
$ws10.Range("R10").Value = 0.0225
$ws10.Range("S10").Value = 0.0075
$ws10.Range("T10").Value = 0.0075
$ws10.Range("U10").Value = 0.0075
$ws10.Range("V10").Value = 0.0225
$ws10.Range("W10").Value = 0.09

# Sheet: Twinsburg Ohio
$ws11 = $wb.Worksheets.Item("Twinsburg Ohio")
$ws11.Range("E2").Value = 0.0585
$ws11.Range("E3").Value = 0.0585
$ws11.Range("E4").Value = 0.0585
$ws11.Range("O4").Value = 0
$ws11.Range("P4").Value = 0
$ws11.Range("Q4").Value = 0
$ws11.Range("R4").Value = 0
$ws11.Range("S4").Value = 0
$ws11.Range("T4").Value = 0
$ws11.Range("U4").Value = 0
$ws11.Range("V4").Value = 0
$ws11.Range("W4").Value = 0
$ws11.Range("O7").Value = ""
$ws11.Range("E8").Value = 0.0969
$ws11.Range("E9").Value = 0.0969
$ws11.Range("E10").Value = 0.0969
$ws11.Range("O10").Value = 0.0769
$ws11.Range("P10").Value = 0.0138416666666667
$ws11.Range("Q10").Value = 0.0138416666666667
$ws11.Range("R10").Value = 0.041525
$ws11.Range("S10").Value = 0.0138416666666667
$ws11.Range("T10").Value = 0.0138416666666667
$ws11.Range("U10").Value = 0.0138416666666667
$ws11.Range("V10").Value = 0.041525
$ws11.Range("W10").Value = 0.1661

# Sheet: Downers Grove Illinois
$ws12 = $wb.Worksheets.Item("Downers Grove Illinois")
$ws12.Range("E2").Value = 0.3448
$ws12.Range("E3").Value = 0.3448
$ws12.Range("E4").Value = 0.3448
$ws12.Range("O4").Value = 0
$ws12.Range("P4").Value = 0
$ws12.Range("Q4").Value = 0
$ws12.Range("R4").Value = 0
$ws12.Range("S4").Value = 0
$ws12.Range("T4").Value = 0
$ws12.Range("U4").Value = 0
$ws12.Range("V4").Value = 0
$ws12.Range("W4").Value = 0
$ws12.Range("O7").Value = ""

Write-Host "Done updating CVDs"